$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 140,141,143,144 (4-cycle swap among matches) ---
# --- and rows 161/162, 169/170 (pairwise swaps), plus row 172 value refresh ---
# Row 140
$ws.Range("B140").Value = 6537915
$ws.Range("F140").Value = "Busan I Park"
$ws.Range("G140").Value = "Chungbuk Cheongju"
$ws.Range("H140").Value = 1
$ws.Range("I140").Value = 1
$ws.Range("J140").Value = "D"
$ws.Range("K140").Value = 1.533
$ws.Range("L140").Value = 4
$ws.Range("M140").Value = 5.25
$ws.Range("N140").Value = 1.444
$ws.Range("O140").Value = 4.2
$ws.Range("P140").Value = 6
$ws.Range("Q140").Value = -1.25
$ws.Range("R140").Value = 1.975
$ws.Range("S140").Value = 1.825
$ws.Range("T140").Value = 2.5
$ws.Range("U140").Value = 1.825
$ws.Range("V140").Value = 1.975
$ws.Range("W140").Value = -1
$ws.Range("X140").Value = 3.2
$ws.Range("Y140").Value = -1
$ws.Range("Z140").Value = -1
$ws.Range("AA140").Value = 0.825
$ws.Range("AB140").Value = -1
$ws.Range("AC140").Value = 0.9750000000000001

# Row 141
$ws.Range("B141").Value = 6537916
$ws.Range("F141").Value = "Gimcheon Sangmu FC"
$ws.Range("G141").Value = "Seoul ELand FC"
$ws.Range("H141").Value = 1
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = "H"
$ws.Range("K141").Value = 1.3
$ws.Range("L141").Value = 5
$ws.Range("M141").Value = 7.5
$ws.Range("N141").Value = 1.222
$ws.Range("O141").Value = 5.5
$ws.Range("P141").Value = 10
$ws.Range("Q141").Value = -1.75
$ws.Range("R141").Value = 1.85
$ws.Range("S141").Value = 1.95
$ws.Range("T141").Value = 3.25
$ws.Range("U141").Value = 1.9
$ws.Range("V141").Value = 1.9
$ws.Range("W141").Value = 0.222
$ws.Range("X141").Value = -1
$ws.Range("Y141").Value = -1
$ws.Range("Z141").Value = -1
$ws.Range("AA141").Value = 0.95
$ws.Range("AB141").Value = -1
$ws.Range("AC141").Value = 0.8999999999999999

# Row 143
$ws.Range("B143").Value = 6527572
$ws.Range("F143").Value = "Bucheon"
$ws.Range("G143").Value = "Jeonnam Dragons"
$ws.Range("H143").Value = 4
$ws.Range("I143").Value = 1
$ws.Range("J143").Value = "H"
$ws.Range("K143").Value = 2.3
$ws.Range("L143").Value = 3.3
$ws.Range("M143").Value = 2.8
$ws.Range("N143").Value = 2.25
$ws.Range("O143").Value = 3.3
$ws.Range("P143").Value = 2.9
$ws.Range("Q143").Value = -0.25
$ws.Range("R143").Value = 1.975
$ws.Range("S143").Value = 1.825
$ws.Range("T143").Value = 2.5
$ws.Range("U143").Value = 1.975
$ws.Range("V143").Value = 1.825
$ws.Range("W143").Value = 1.25
$ws.Range("X143").Value = -1
$ws.Range("Y143").Value = -1
$ws.Range("Z143").Value = 0.9750000000000001
$ws.Range("AA143").Value = -1
$ws.Range("AB143").Value = 0.9750000000000001
$ws.Range("AC143").Value = -1

# Row 144
$ws.Range("B144").Value = 6531883
$ws.Range("F144").Value = "Seongnam FC"
$ws.Range("G144").Value = "Ansan Greeners FC"
$ws.Range("H144").Value = 0
$ws.Range("I144").Value = 2
$ws.Range("J144").Value = "A"
$ws.Range("K144").Value = 1.7
$ws.Range("L144").Value = 3.8
$ws.Range("M144").Value = 4.2
$ws.Range("N144").Value = 1.75
$ws.Range("O144").Value = 4
$ws.Range("P144").Value = 3.75
$ws.Range("Q144").Value = -0.75
$ws.Range("R144").Value = 2
$ws.Range("S144").Value = 1.8
$ws.Range("T144").Value = 3.25
$ws.Range("U144").Value = 2
$ws.Range("V144").Value = 1.8
$ws.Range("W144").Value = -1
$ws.Range("X144").Value = -1
$ws.Range("Y144").Value = 2.75
$ws.Range("Z144").Value = -1
$ws.Range("AA144").Value = 0.8
$ws.Range("AB144").Value = -1
$ws.Range("AC144").Value = 0.8

# Row 161
$ws.Range("B161").Value = 7738657
$ws.Range("F161").Value = "Bucheon"
$ws.Range("G161").Value = "Seoul ELand FC"
$ws.Range("H161").Value = 1
$ws.Range("I161").Value = 0
$ws.Range("J161").Value = "H"
$ws.Range("K161").Value = 2.5
$ws.Range("L161").Value = 3.2
$ws.Range("M161").Value = 2.5
$ws.Range("N161").Value = 2.45
$ws.Range("O161").Value = 3.2
$ws.Range("P161").Value = 2.55
$ws.Range("Q161").Value = 0
$ws.Range("R161").Value = 1.875
$ws.Range("S161").Value = 1.925
$ws.Range("T161").Value = 2.25
$ws.Range("U161").Value = 1.825
$ws.Range("V161").Value = 1.975
$ws.Range("W161").Value = 1.45
$ws.Range("X161").Value = -1
$ws.Range("Y161").Value = -1
$ws.Range("Z161").Value = 0.875
$ws.Range("AA161").Value = -1
$ws.Range("AB161").Value = -1
$ws.Range("AC161").Value = 0.9750000000000001

# Row 162
$ws.Range("B162").Value = 7738682
$ws.Range("F162").Value = "Gyeongnam FC"
$ws.Range("G162").Value = "Chungnam Asan FC"
$ws.Range("H162").Value = 1
$ws.Range("I162").Value = 2
$ws.Range("J162").Value = "A"
$ws.Range("K162").Value = 2.25
$ws.Range("L162").Value = 3.25
$ws.Range("M162").Value = 2.75
$ws.Range("N162").Value = 2
$ws.Range("O162").Value = 3.3
$ws.Range("P162").Value = 3.2
$ws.Range("Q162").Value = -0.25
$ws.Range("R162").Value = 1.8
$ws.Range("S162").Value = 2
$ws.Range("T162").Value = 2.25
$ws.Range("U162").Value = 1.875
$ws.Range("V162").Value = 1.925
$ws.Range("W162").Value = -1
$ws.Range("X162").Value = -1
$ws.Range("Y162").Value = 2.2
$ws.Range("Z162").Value = -1
$ws.Range("AA162").Value = 1
$ws.Range("AB162").Value = 0.875
$ws.Range("AC162").Value = -1

# Row 169
$ws.Range("B169").Value = 7738685
$ws.Range("F169").Value = "Bucheon"
$ws.Range("G169").Value = "Gyeongnam FC"
$ws.Range("H169").Value = 2
$ws.Range("I169").Value = 0
$ws.Range("J169").Value = "H"
$ws.Range("K169").Value = 2.5
$ws.Range("L169").Value = 3.1
$ws.Range("M169").Value = 2.55
$ws.Range("N169").Value = 2.55
$ws.Range("O169").Value = 3.2
$ws.Range("P169").Value = 2.6
$ws.Range("Q169").Value = 0
$ws.Range("R169").Value = 1.9
$ws.Range("S169").Value = 1.9
$ws.Range("T169").Value = 2.25
$ws.Range("U169").Value = 1.975
$ws.Range("V169").Value = 1.825
$ws.Range("W169").Value = 1.55
$ws.Range("X169").Value = -1
$ws.Range("Y169").Value = -1
$ws.Range("Z169").Value = 0.8999999999999999
$ws.Range("AA169").Value = -1
$ws.Range("AB169").Value = -0.5
$ws.Range("AC169").Value = 0.4125

# Row 170
$ws.Range("B170").Value = 7738659
$ws.Range("F170").Value = "Jeonnam Dragons"
$ws.Range("G170").Value = "Cheonan City"
$ws.Range("H170").Value = 1
$ws.Range("I170").Value = 0
$ws.Range("J170").Value = "H"
$ws.Range("K170").Value = 1.85
$ws.Range("L170").Value = 3.5
$ws.Range("M170").Value = 3.5
$ws.Range("N170").Value = 1.75
$ws.Range("O170").Value = 3.6
$ws.Range("P170").Value = 3.75
$ws.Range("Q170").Value = -0.5
$ws.Range("R170").Value = 1.825
$ws.Range("S170").Value = 1.975
$ws.Range("T170").Value = 2.5
$ws.Range("U170").Value = 2
$ws.Range("V170").Value = 1.8
$ws.Range("W170").Value = 0.75
$ws.Range("X170").Value = -1
$ws.Range("Y170").Value = -1
$ws.Range("Z170").Value = 0.825
$ws.Range("AA170").Value = -1
$ws.Range("AB170").Value = -1
$ws.Range("AC170").Value = 0.8

# Row 172
$ws.Range("B172").Value = 7738686
$ws.Range("F172").Value = "FC Anyang"
$ws.Range("G172").Value = "Ansan Greeners FC"
$ws.Range("K172").Value = 1.55
$ws.Range("L172").Value = 4
$ws.Range("M172").Value = 5.5
$ws.Range("N172").Value = 1.5
$ws.Range("O172").Value = 4.2
$ws.Range("P172").Value = 6
$ws.Range("Q172").Value = -1
$ws.Range("R172").Value = 1.775
$ws.Range("S172").Value = 2.025
$ws.Range("T172").Value = 2.75
$ws.Range("U172").Value = 1.85
$ws.Range("V172").Value = 1.95
$ws.Range("W172").Value = 0
$ws.Range("X172").Value = 0
$ws.Range("Y172").Value = 0
$ws.Range("Z172").Value = 0
$ws.Range("AA172").Value = 0

# --- Append new row 173 ---
$ws.Range("A172").Copy()
$ws.Range("A173").PasteSpecial(-4122)
$ws.Range("E172").Copy()
$ws.Range("E173").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A173").Value = 171
$ws.Range("B173").Value = 7738661
$ws.Range("C173").Value = "South Korea K League 2"
$ws.Range("D173").Value = "South Korea K League 2"
$ws.Range("E173").Value = 45389.08333333334
$ws.Range("F173").Value = "Cheonan City"
$ws.Range("G173").Value = "Gyeongnam FC"
$ws.Range("K173").Value = 4
$ws.Range("L173").Value = 3.4
$ws.Range("M173").Value = 1.833
$ws.Range("N173").Value = 3.75
$ws.Range("O173").Value = 3.3
$ws.Range("P173").Value = 1.909
$ws.Range("Q173").Value = 0.5
$ws.Range("R173").Value = 1.825
$ws.Range("S173").Value = 1.975
$ws.Range("T173").Value = 2.5
$ws.Range("U173").Value = 1.975
$ws.Range("V173").Value = 1.825
$ws.Range("W173").Value = 0
$ws.Range("X173").Value = 0
$ws.Range("Y173").Value = 0
$ws.Range("Z173").Value = 0
$ws.Range("AA173").Value = 0
